$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "C003"
$ws.Range("B4").Value = "Abdel"
$ws.Range("C4").Value = 4582258
$ws.Range("D4").Value = 779525255
